# Apply the BSR bill updates to rows 8-14 (item detail rows) and the
# grand-total rows 16 & 18, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay a TEXT cell even when it looks like
# a number (e.g. "2", "3.0", "11264.00"). Forcing the NumberFormat to "@"
# (Text) before assigning keeps Excel from silently re-typing the cell as
# a numeric value.
function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Helper: write an ordinary string (never ambiguous with a number) - plain
# assignment is fine here.
function Set-StringValue($rangeAddress, $value) {
    $ws.Range($rangeAddress).Value = $value
}

# Helper: write a real number.
function Set-NumberValue($rangeAddress, $value) {
    $ws.Range($rangeAddress).Value = $value
}

# ---------------------------------------------------------------------
# Row 8
# ---------------------------------------------------------------------
Set-StringValue "A8" "P. point"
Set-NumberValue "C8" 44
Set-TextValue   "D8" "2"
Set-StringValue "E8" "Short point (up to 3 mtr.)"
Set-NumberValue "F8" 256
Set-TextValue   "G8" "11264.00"

# ---------------------------------------------------------------------
# Row 9
# ---------------------------------------------------------------------
Set-NumberValue "C9" 21
Set-TextValue   "D9" "3"
Set-StringValue "E9" "Medium point (up to 6 mtr.)"
Set-NumberValue "F9" 472
Set-TextValue   "G9" "9912.00"

# ---------------------------------------------------------------------
# Row 10
# ---------------------------------------------------------------------
Set-NumberValue "C10" 78
Set-TextValue   "G10" "51636.00"

# ---------------------------------------------------------------------
# Row 11
# ---------------------------------------------------------------------
Set-NumberValue "C11" 87

# ---------------------------------------------------------------------
# Row 12
# ---------------------------------------------------------------------
Set-StringValue "A12" "Each"
Set-NumberValue "C12" 16
Set-TextValue   "D12" "3.0"
Set-StringValue "E12" 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumberValue "F12" 23
Set-TextValue   "G12" "368.00"

# ---------------------------------------------------------------------
# Row 13
# ---------------------------------------------------------------------
Set-StringValue "A13" ""
Set-NumberValue "C13" 18
Set-TextValue   "D13" "8"
Set-StringValue "E13" "Total"

# ---------------------------------------------------------------------
# Row 14
# ---------------------------------------------------------------------
Set-StringValue "A14" "%"
Set-NumberValue "C14" 26
Set-TextValue   "D14" "9"
Set-StringValue "E14" "Add Tender Premium "

# ---------------------------------------------------------------------
# Grand totals (rows 16 & 18)
# ---------------------------------------------------------------------
Set-TextValue "G16" "73180.00"
Set-TextValue "H16" "73180.00"
Set-TextValue "G18" "73180.00"
Set-TextValue "H18" "73180.00"
